# RoomNumbers.xlsx edit:
#  - Header rename: "Seating Capacity" -> "Capacity", "Biometric" -> "Biometric?"
#  - Header row (A1:C1) made bold; B1/C1 centered
#  - Numeric/flag columns (B, C) center-aligned for all data rows
#  - New/explicit width for column C
#  - Room D217 (150 seats, biometric) moved up so it sits right after D216,
#    with D218/D220/D222/D226 shifting down one row (consecutive classes
#    now get the same room block)
#  - View reset to top of sheet, selection on C1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header text ---
$ws.Range("B1").Value = "Capacity"
$ws.Range("C1").Value = "Biometric?"

# --- Header styling: bold for A1:C1, centered for B1:C1 ---
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("B1:C1").HorizontalAlignment = -4108

# --- Center-align the data in columns B and C ---
$ws.Columns.Item(2).HorizontalAlignment = -4108
$ws.Columns.Item(3).HorizontalAlignment = -4108

# --- Column C width (new dedicated width, close to the 15.21 target -
#     Excel's ColumnWidth is quantized to pixel steps so this lands on the
#     nearest achievable width) ---
$ws.Columns.Item(3).ColumnWidth = 14.3

# --- Reorder rows 34-38: pull D217 (150, biometric) up so it follows D216,
#     push D218/D220/D222/D226 down one row each ---
$ws.Range("A34").Value = "D217"
$ws.Range("B34").Value = 150
$ws.Range("C34").Value = 1

$ws.Range("A35").Value = "D218"
$ws.Range("B35").Value = 30
$ws.Range("C35").Value = 0

$ws.Range("A36").Value = "D220"
$ws.Range("B36").Value = 30
$ws.Range("C36").Value = 0

$ws.Range("A37").Value = "D222"
$ws.Range("B37").Value = 60
$ws.Range("C37").Value = 0

$ws.Range("A38").Value = "D226"
$ws.Range("B38").Value = 60
$ws.Range("C38").Value = 0

# --- View: scroll back to top, select C1 ---
$ws.Application.Goto($ws.Range("A1")) | Out-Null
$ws.Range("C1").Select() | Out-Null
